# Update the "Metadata" worksheet (sheet1 in the OOXML package):
#  - bump the Version value
#  - bump the Date value
#  - insert a new "Jurisdiction" property row (with an empty value) right
#    before the "Description" row, shifting everything below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version (row 3, column B)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date (row 8, column B)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row above the existing "Description" row (row 11) so that the
# following rows (Description, Purpose, Copyright, ...) shift down by one.
$ws.Rows.Item(11).Insert()

# Carry over the formatting used by the other property rows (row 10, "Contact")
# onto the freshly inserted (otherwise unformatted) row.
$ws.Range("A10:B10").Copy($ws.Range("A11:B11"))

# Populate the newly inserted row 11 with the Jurisdiction property.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
